$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F4 369->370, F5 417->418, F9 6284->6291
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 370
$ws1.Range("F5").Value = 418
$ws1.Range("F9").Value = 6291

# Sheet "全部类型" (sheet4): F4 369->370, F5 417->418, F11 6284->6291
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 370
$ws4.Range("F5").Value = 418
$ws4.Range("F11").Value = 6291
